# This workbook's price/volume rows (D:H) and ticker column (I) for most
# rows were populated from the wrong source company ("extra files" bug per
# the commit message). Every row except row 11 (already correct) must be
# corrected to the true UBER OHLC/shares data, and the ticker column reset
# to "UBER". Once nothing references the stray ticker strings (FTNT, TSM,
# VRNT, ...) anymore, the shared-string table naturally shrinks back down
# on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => open, close, high, low, shares_outstanding
$fixedRows = [ordered]@{
  2  = @(42,                  40.40999984741211,  45,                 36.08000183105469,  2085418676)
  3  = @(42,                  40.40999984741211,  45,                 36.08000183105469,  2085418676)
  4  = @(42,                  40.40999984741211,  45,                 36.08000183105469,  2085418676)
  5  = @(42,                  40.40999984741211,  45,                 36.08000183105469,  2085418676)
  6  = @(42,                  40.40999984741211,  45,                 36.08000183105469,  2085418676)
  7  = @(42,                  40.40999984741211,  45,                 36.08000183105469,  2085418676)
  8  = @(46.97999954223633,   42.13999938964844,  47.04000091552734,  41.84000015258789,  2085418676)
  9  = @(30.3700008392334,    31.5,               33.88999938964844,  28.30999946594238,  2085418676)
  10 = @(29.94000053405762,   36.29000091552734,  37.95000076293945,  29.79000091552734,  2085418676)
  12 = @(30.95999908447266,   30.26000022888184,  34.45999908447266,  29.78499984741211,  2085418676)
  13 = @(36.86000061035156,   33.40999984741211,  38.34999847412109,  32.90000152587891,  2085418676)
  14 = @(52.22000122070312,   50.93000030517578,  60.02999877929688,  47.15000152587891,  2085418676)
  15 = @(55.61000061035156,   54.77000045776367,  61.5,               53.11000061035156,  2085418676)
  16 = @(50.88000106811523,   43.45999908447266,  52.02999877929688,  43.34000015258789,  2085418676)
  17 = @(45.91500091552734,   43.81999969482422,  48.88000106811523,  43.56999969482422,  2085418676)
  18 = @(42.47999954223633,   37.40000152587891,  45.90000152587891,  32.81000137329102,  2085418676)
  19 = @(35.68000030517578,   31.47999954223633,  36.93000030517578,  30.15999984741211,  2085418676)
  20 = @(20.70000076293945,   23.45000076293945,  24.80999946594238,  20.43000030517578,  2085418676)
  21 = @(26.76000022888184,   26.56999969482422,  30.23999977111816,  22.94000053405762,  2085418676)
  22 = @(25.3700008392334,    30.93000030517578,  31.20000076293945,  25.02000045776367,  2085418676)
  23 = @(31.5,                31.04999923706055,  32.93000030517578,  29.21999931335449,  2085418676)
  24 = @(43.16999816894531,   49.45999908447266,  49.4900016784668,   41.74599838256836,  2085418676)
  25 = @(45.56999969482422,   43.27999877929688,  47.27500152587891,  40.09000015258789,  2085418676)
  26 = @(60.72999954223633,   65.26999664306641,  67.69000244140625,  57.22000122070312,  2085418676)
  27 = @(77,                  66.26999664306641,  78.84999847412109,  66.06999969482422,  2085418676)
  28 = @(72.19999694824219,   64.47000122070312,  75.40000152587891,  62.61000061035156,  2085418676)
  29 = @(75.76999664306641,   72.05000305175781,  87,                 69.87000274658203,  2085418676)
  30 = @(62.18999862670898,   66.84999847412109,  69.67500305175781,  61.31000137329102,  2085418676)
  31 = @(72.59999847412109,   81.01000213623047,  81.21900177001953,  60.63000106811523,  2085418676)
  32 = @(92.84999847412109,   87.75,              97.71499633789062,  86.36000061035156,  2085418676)
}

foreach ($row in $fixedRows.Keys) {
    $vals = $fixedRows[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
    $ws.Range("H$row").Value = $vals[4]
    $ws.Range("I$row").Value = "UBER"
}
